# Auto-generated Excel COM-interop script that replays a scheduled-runner
# market-price refresh across the workbook's 8 job-sheets (ALC, ARM, BSM, CRP,
# CUL, GSM, LTW, WVR). Each row holds "currentAveragePrice*" / "LevePrice*" /
# "LeveProfit*" figures (columns H:N) pulled from an external market-data feed;
# this script rewrites the refreshed numbers cell-by-cell per row, matching the
# upstream data source's snapshot. A few rows also gain/lose a cell (the feed
# only emits a LeveProfitNQ or LeveProfitHQ value when that variant is craftable).

$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")

# Row 8
$ws.Range("H8").Value = 3201.25
$ws.Range("I8").Value = 55
$ws.Range("K8").Value = 165
$ws.Range("M8").Value = -26

# Row 18
$ws.Range("H18").Value = 2044.5
$ws.Range("I18").Value = 2059.3333
$ws.Range("K18").Value = 2059.3333
$ws.Range("M18").Value = -1775.3333

# Row 40
$ws.Range("H40").Value = 5240.3706

# Row 42
$ws.Range("H42").Value = 40.6
$ws.Range("I42").Value = 46.833332
$ws.Range("J42").Value = 31.25
$ws.Range("K42").Value = 140.499996
$ws.Range("L42").Value = 93.75
$ws.Range("M42").Value = 89.50000399999999
$ws.Range("N42").Value = -553.75

# Row 113
$ws.Range("H113").Value = 3749.625
$ws.Range("I113").Value = 3666.1667
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 3666.1667
$ws.Range("L113").Value = 4000
$ws.Range("M113").Value = -412.1667000000002
$ws.Range("N113").Value = -10508

# Row 137
$ws.Range("H137").Value = 33342742
$ws.Range("I137").Value = 41669252
$ws.Range("J137").Value = 36700
$ws.Range("K137").Value = 125007756
$ws.Range("L137").Value = 110100
$ws.Range("M137").Value = -125005206
$ws.Range("N137").Value = -115200

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")

# Row 28
$ws.Range("H28").Value = 25000.25
$ws.Range("I28").Value = 20000
$ws.Range("K28").Value = 20000
$ws.Range("M28").Value = -19808

# Row 32
$ws.Range("H32").Value = 913569.6
$ws.Range("I32").Value = 975800.9
$ws.Range("K32").Value = 975800.9
$ws.Range("M32").Value = -975513.9

# Row 99
$ws.Range("H99").Value = 25000.25
$ws.Range("I99").Value = 20000
$ws.Range("K99").Value = 20000
$ws.Range("M99").Value = -17005

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")

# Row 20
$ws.Range("H20").Value = 26067.678
$ws.Range("I20").Value = 29085.898
$ws.Range("J20").Value = 14691.308
$ws.Range("K20").Value = 29085.898
$ws.Range("L20").Value = 14691.308
$ws.Range("M20").Value = -28838.898
$ws.Range("N20").Value = -15185.308

# Row 99
$ws.Range("H99").Value = 5974.75
$ws.Range("I99").Value = 6879.8237
$ws.Range("J99").Value = 846
$ws.Range("K99").Value = 6879.8237
$ws.Range("L99").Value = 846
$ws.Range("M99").Value = -5381.8237
$ws.Range("N99").Value = -3842

# Row 105
$ws.Range("H105").Value = 8679.799999999999
$ws.Range("I105").Value = 2500
$ws.Range("J105").Value = 17949.5
$ws.Range("K105").Value = 2500
$ws.Range("L105").Value = 17949.5
$ws.Range("M105").Value = -753
$ws.Range("N105").Value = -21443.5

# Row 134
$ws.Range("H134").Value = 2033649.6
$ws.Range("I134").Value = 1169.8462
$ws.Range("K134").Value = 3509.5386
$ws.Range("M134").Value = -974.5385999999999

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")

# Row 7
$ws.Range("H7").Value = 144.94444
$ws.Range("I7").Value = 123.833336
$ws.Range("J7").Value = 187.16667
$ws.Range("K7").Value = 123.833336
$ws.Range("L7").Value = 187.16667
$ws.Range("M7").Value = -10.833336
$ws.Range("N7").Value = -413.16667

# Row 31
$ws.Range("H31").Value = 3676686.5
$ws.Range("I31").Value = 1468263.8
$ws.Range("K31").Value = 1468263.8
$ws.Range("M31").Value = -1467968.8

# Row 33
$ws.Range("H33").Value = 2599
$ws.Range("I33").Value = 2599
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 2599
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -2220
$ws.Range("N33").ClearContents()

# Row 34
$ws.Range("H34").Value = 3676686.5
$ws.Range("I34").Value = 1468263.8
$ws.Range("K34").Value = 1468263.8
$ws.Range("M34").Value = -1468061.8

# Row 35
$ws.Range("H35").Value = 1985
$ws.Range("I35").Value = 2199.8
$ws.Range("J35").Value = 1448
$ws.Range("K35").Value = 2199.8
$ws.Range("L35").Value = 1448
$ws.Range("M35").Value = -1905.8
$ws.Range("N35").Value = -2036

# Row 38
$ws.Range("H38").Value = 3863.6365
$ws.Range("J38").Value = 3863.6365
$ws.Range("L38").Value = 3863.6365
$ws.Range("N38").Value = -4617.636500000001

# Row 39
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()

# Row 46
$ws.Range("H46").Value = 3863.6365
$ws.Range("J46").Value = 3863.6365
$ws.Range("L46").Value = 3863.6365
$ws.Range("N46").Value = -4285.636500000001

# Row 49
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()

# Row 58
$ws.Range("H58").Value = 17350110
$ws.Range("I58").Value = 23813270
$ws.Range("K58").Value = 23813270
$ws.Range("M58").Value = -23813067

# Row 132
$ws.Range("H132").Value = 3309.077
$ws.Range("I132").Value = 3020.875
$ws.Range("J132").Value = 3770.2
$ws.Range("K132").Value = 9062.625
$ws.Range("L132").Value = 11310.6
$ws.Range("M132").Value = -6532.625
$ws.Range("N132").Value = -16370.6

# Row 136
$ws.Range("H136").Value = 17350110
$ws.Range("I136").Value = 23813270
$ws.Range("K136").Value = 71439810
$ws.Range("M136").Value = -71437260

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")

# Row 5
$ws.Range("H5").Value = 2359504.2
$ws.Range("I5").Value = 1786214.6
$ws.Range("J5").Value = 4652662
$ws.Range("K5").Value = 5358643.800000001
$ws.Range("L5").Value = 13957986
$ws.Range("M5").Value = -5358531.800000001
$ws.Range("N5").Value = -13958210

# Row 124
$ws.Range("H124").Value = 10925.637
$ws.Range("I124").Value = 7530.3335
$ws.Range("K124").Value = 22591.0005
$ws.Range("M124").Value = -17681.0005

# Row 131
$ws.Range("H131").Value = 4347.653
$ws.Range("I131").Value = 550.8461
$ws.Range("J131").Value = 5718.722
$ws.Range("K131").Value = 1652.5383
$ws.Range("L131").Value = 17156.166
$ws.Range("M131").Value = 3387.4617
$ws.Range("N131").Value = -27236.166

# Row 135
$ws.Range("H135").Value = 2359504.2
$ws.Range("I135").Value = 1786214.6
$ws.Range("J135").Value = 4652662
$ws.Range("K135").Value = 16075931.4
$ws.Range("L135").Value = 41873958
$ws.Range("M135").Value = -16073396.4
$ws.Range("N135").Value = -41879028

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")

# Row 31
$ws.Range("H31").Value = 1916.5
$ws.Range("I31").Value = 1916.5
$ws.Range("K31").Value = 1916.5
$ws.Range("M31").Value = -1624.5

# Row 37
$ws.Range("H37").Value = 1916.5
$ws.Range("I37").Value = 1916.5
$ws.Range("K37").Value = 1916.5
$ws.Range("M37").Value = -1639.5

# Row 68
$ws.Range("H68").Value = 34000
$ws.Range("J68").Value = 34000
$ws.Range("L68").Value = 34000
$ws.Range("N68").Value = -35622

# Row 71
$ws.Range("H71").Value = 34000
$ws.Range("J71").Value = 34000
$ws.Range("L71").Value = 102000
$ws.Range("N71").Value = -110112

# Row 102
$ws.Range("H102").Value = 1989.7778
$ws.Range("I102").Value = 1986.3334
$ws.Range("J102").Value = 1996.6666
$ws.Range("K102").Value = 1986.3334
$ws.Range("L102").Value = 1996.6666
$ws.Range("M102").Value = -364.3334
$ws.Range("N102").Value = -5240.6666

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")

# Row 16
$ws.Range("H16").Value = 787.5
$ws.Range("I16").Value = 726.4211
$ws.Range("K16").Value = 726.4211
$ws.Range("M16").Value = -556.4211

# Row 22
$ws.Range("H22").Value = 2667
$ws.Range("I22").Value = 2229.9285
$ws.Range("J22").Value = 3026.9412
$ws.Range("K22").Value = 2229.9285
$ws.Range("L22").Value = 3026.9412
$ws.Range("M22").Value = -1934.9285
$ws.Range("N22").Value = -3616.9412

# Row 27
$ws.Range("H27").Value = 2667
$ws.Range("I27").Value = 2229.9285
$ws.Range("J27").Value = 3026.9412
$ws.Range("K27").Value = 2229.9285
$ws.Range("L27").Value = 3026.9412
$ws.Range("M27").Value = -2122.9285
$ws.Range("N27").Value = -3240.9412

# Row 46
$ws.Range("H46").Value = 4247.15
$ws.Range("J46").Value = 5499.2144
$ws.Range("L46").Value = 5499.2144
$ws.Range("N46").Value = -5875.2144

# Row 55
$ws.Range("H55").Value = 1600.5333
$ws.Range("I55").Value = 1541.3077
$ws.Range("J55").Value = 1645.8235
$ws.Range("K55").Value = 1541.3077
$ws.Range("L55").Value = 1645.8235
$ws.Range("M55").Value = -1368.3077
$ws.Range("N55").Value = -1991.8235

# Row 82
$ws.Range("H82").Value = 3950
$ws.Range("I82").Value = 5266.6665
$ws.Range("J82").Value = 2962.5
$ws.Range("K82").Value = 5266.6665
$ws.Range("L82").Value = 2962.5
$ws.Range("M82").Value = -4905.6665
$ws.Range("N82").Value = -3684.5

# Row 85
$ws.Range("H85").Value = 3950
$ws.Range("I85").Value = 5266.6665
$ws.Range("J85").Value = 2962.5
$ws.Range("K85").Value = 5266.6665
$ws.Range("L85").Value = 2962.5
$ws.Range("M85").Value = -4018.6665
$ws.Range("N85").Value = -5458.5

# Row 122
$ws.Range("H122").Value = 3446.35
$ws.Range("I122").Value = 3067.5715
$ws.Range("K122").Value = 9202.7145
$ws.Range("M122").Value = -6752.7145

# Row 136
$ws.Range("H136").Value = 10872991
$ws.Range("I136").Value = 6582231
$ws.Range("J136").Value = 31254102
$ws.Range("K136").Value = 19746693
$ws.Range("L136").Value = 93762306
$ws.Range("M136").Value = -19744143
$ws.Range("N136").Value = -93767406

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")

# Row 81
$ws.Range("H81").Value = 130394.5
$ws.Range("J81").Value = 171776
$ws.Range("L81").Value = 343552
$ws.Range("N81").Value = -345674

# Row 84
$ws.Range("H84").Value = 130394.5
$ws.Range("J84").Value = 171776
$ws.Range("L84").Value = 1717760
$ws.Range("N84").Value = -1728368

# Row 122
$ws.Range("H122").Value = 68749.53
$ws.Range("I122").Value = 2785.9167
$ws.Range("K122").Value = 8357.750100000001
$ws.Range("M122").Value = -5907.750100000001
